$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = -1
$ws.Range("F15").Value = 2
$ws.Range("F18").Value = -4
$ws.Range("F26").Value = -3
$ws.Range("F29").Value = 0
